$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a shared-string "text" value (e.g. the literal text
# "0" or "***.*" used as a placeholder in this report) while preserving the
# destination cell's existing number format/style. Directly assigning a
# numeric-looking string via .Value gets auto-coerced to a real number by
# Excel, so instead we copy a donor cell that already holds the desired text
# (with the right style) and paste format + value separately.
function Set-TextCell {
    param($targetAddr, $sourceAddr)
    $ws.Range($sourceAddr).Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range($sourceAddr).Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# ---------------------------------------------------------------------
# Header text updates (Volume/Number and report week dates)
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  9"
$ws.Range("C9").Value = "Report Covering the Week  2/27/2023  Through  3/5/2023"

# ---------------------------------------------------------------------
# Crime-complaints table: refreshed weekly figures
# ---------------------------------------------------------------------
$ws.Range("F14").Value = 1
$ws.Range("L15").Value = -50
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = -75
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -46.153846153846
$ws.Range("I16").Value = 22
$ws.Range("J16").Value = 22
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 144.444444444444
$ws.Range("M16").Value = -35.294117647058
$ws.Range("N16").Value = -78.217821782178
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 25
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = -22.727272727272
$ws.Range("I17").Value = 37
$ws.Range("J17").Value = 52
$ws.Range("K17").Value = -28.846153846153
$ws.Range("L17").Value = 37.037037037037
$ws.Range("M17").Value = 54.166666666666
$ws.Range("N17").Value = 32.142857142857
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = -12.5
$ws.Range("I18").Value = 11
$ws.Range("J18").Value = 13
$ws.Range("K18").Value = -15.384615384615
$ws.Range("L18").Value = -8.333333333333
$ws.Range("M18").Value = -75.555555555555
$ws.Range("N18").Value = -90.598290598290
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 11.111111111111
$ws.Range("F19").Value = 19
$ws.Range("G19").Value = 24
$ws.Range("H19").Value = -20.833333333333
$ws.Range("I19").Value = 44
$ws.Range("J19").Value = 43
$ws.Range("K19").Value = 2.325581395348
$ws.Range("L19").Value = 109.52380952381
$ws.Range("M19").Value = 22.222222222222
$ws.Range("N19").Value = -29.032258064516
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 400
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = 6.666666666666
$ws.Range("I20").Value = 26
$ws.Range("J20").Value = 30
$ws.Range("K20").Value = -13.333333333333
$ws.Range("L20").Value = 100
$ws.Range("M20").Value = 23.809523809523
$ws.Range("N20").Value = -93.467336683417
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = 25
$ws.Range("F21").Value = 67
$ws.Range("G21").Value = 83
$ws.Range("H21").Value = -19.277108433734
$ws.Range("I21").Value = 143
$ws.Range("J21").Value = 164
$ws.Range("K21").Value = -12.804878048780
$ws.Range("L21").Value = 68.235294117647
$ws.Range("M21").Value = -11.180124223602
$ws.Range("N21").Value = -79.859154929577
$ws.Range("C23").Value = 4
$ws.Range("F23").Value = 10
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = 66.666666666666
$ws.Range("I23").Value = 18
$ws.Range("J23").Value = 11
$ws.Range("K23").Value = 63.636363636363
$ws.Range("L23").Value = 38.461538461538
$ws.Range("M23").Value = 500
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = -7.692307692307
$ws.Range("F24").Value = 52
$ws.Range("G24").Value = 43
$ws.Range("H24").Value = 20.930232558139
$ws.Range("I24").Value = 122
$ws.Range("J24").Value = 98
$ws.Range("K24").Value = 24.489795918367
$ws.Range("L24").Value = 50.617283950617
$ws.Range("M24").Value = 74.285714285714
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 19
$ws.Range("G25").Value = 18
$ws.Range("H25").Value = 5.555555555555
$ws.Range("I25").Value = 49
$ws.Range("J25").Value = 51
$ws.Range("K25").Value = -3.921568627450
$ws.Range("L25").Value = 53.125
$ws.Range("M25").Value = -16.949152542372
$ws.Range("G26").Value = 3
$ws.Range("J26").Value = 4
$ws.Range("K26").Value = -75
$ws.Range("L26").Value = -75
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 2
$ws.Range("J27").Value = 6
$ws.Range("K27").Value = -66.666666666666
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 100
$ws.Range("L28").Value = -20
$ws.Range("M28").Value = 33.333333333333
$ws.Range("N28").Value = -20
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 100
$ws.Range("L29").Value = -20
$ws.Range("M29").Value = 100
$ws.Range("N29").Value = -20

# Cells that switch to the "{0 / ***.* }" text placeholders (copied from
# stable donor cells C14 "0" and E14 "***.*", both style 14)
Set-TextCell "G14" "C14"
Set-TextCell "H14" "E14"
Set-TextCell "C22" "C14"
Set-TextCell "D22" "C14"
Set-TextCell "E22" "E14"
Set-TextCell "D23" "C14"
Set-TextCell "E23" "E14"
Set-TextCell "C26" "C14"
Set-TextCell "D26" "C14"
Set-TextCell "E26" "E14"
Set-TextCell "F26" "C14"
Set-TextCell "C27" "C14"
Set-TextCell "D27" "C14"
Set-TextCell "E27" "E14"

$ws.Application.CutCopyMode = $false
